# Scenarios.xlsx edit: add a new individual ("Indiv1") to the
# PopulationScenario / PopulationScenarioFromCSV test rows so the
# populationTimeProfile plotting tests have data to exercise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Rows 4 and 5 (PopulationScenario / PopulationScenarioFromCSV) used the
# "Indiv" individual in column B (IndividualId) -- point them at "Indiv1"
# instead, matching the individual already used in row 2.
$ws.Range("B4").Value = "Indiv1"
$ws.Range("B5").Value = "Indiv1"

# Leave the cursor where the author left it when saving.
$ws.Activate()
$ws.Range("E19").Select() | Out-Null
